# Applies the "added equations to get population weighted acmr" edit:
#  - Adds a column of China (sub-national) location names next to the
#    location_id column on the "location_id" sheet.
#  - Clears out the stray formatted-but-empty cells that used to trail
#    row 4 (J4:N4), replacing them with a smaller formatted block (D4:I4).
#  - Updates sheet selections / active-tab bookkeeping to match the
#    state the workbook was left in after the edit.

$wb = $excel.ActiveWorkbook

$wsAge = $wb.Worksheets.Item("age_groups")
$wsLoc = $wb.Worksheets.Item("location_id")

# --- location_id sheet: add province / SAR names in column C -------------
# Filled in the same order the author originally typed them in (matches
# the resulting shared-string table ordering): the alphabetic run of
# mainland provinces first, then "Inner Mongolia" (out of alphabetic
# order, location_id 502) and finally "Hong Kong" (location_id 354) at
# the very end.
$wsLoc.Range("C4").Value = "Macao"
$wsLoc.Range("C5").Value = "Anhui"
$wsLoc.Range("C6").Value = "Beijing"
$wsLoc.Range("C7").Value = "Chongqing"
$wsLoc.Range("C8").Value = "Fujian"
$wsLoc.Range("C9").Value = "Gansu"
$wsLoc.Range("C10").Value = "Guangdong"
$wsLoc.Range("C11").Value = "Guangxi"
$wsLoc.Range("C12").Value = "Guizhou"
$wsLoc.Range("C13").Value = "Hainan"
$wsLoc.Range("C14").Value = "Hebei"
$wsLoc.Range("C15").Value = "Heilongjiang"
$wsLoc.Range("C17").Value = "Henan"
$wsLoc.Range("C18").Value = "Hubei"
$wsLoc.Range("C19").Value = "Hunan"
$wsLoc.Range("C20").Value = "Jiangsu"
$wsLoc.Range("C21").Value = "Jiangxi"
$wsLoc.Range("C22").Value = "Jilin"
$wsLoc.Range("C23").Value = "Liaoning"
$wsLoc.Range("C24").Value = "Ningxia"
$wsLoc.Range("C25").Value = "Qinghai"
$wsLoc.Range("C26").Value = "Shaanxi"
$wsLoc.Range("C27").Value = "Shandong"
$wsLoc.Range("C28").Value = "Shanghai"
$wsLoc.Range("C29").Value = "Shanxi"
$wsLoc.Range("C30").Value = "Sichuan"
$wsLoc.Range("C31").Value = "Tianjin"
$wsLoc.Range("C32").Value = "Tibet"
$wsLoc.Range("C33").Value = "Xinjiang"
$wsLoc.Range("C34").Value = "Yunnan"
$wsLoc.Range("C35").Value = "Zhejiang"
$wsLoc.Range("C16").Value = "Inner Mongolia"
$wsLoc.Range("C3").Value = "Hong Kong"

# --- location_id sheet: tidy up row 4's trailing empty formatted cells ---
# Originally I4:N4 held a centered blank style; now only D4:I4 should.
$wsLoc.Range("J4:N4").Clear() | Out-Null
$wsLoc.Range("D4:I4").HorizontalAlignment = -4108

# --- column widths ---------------------------------------------------
$wsLoc.Columns.Item(3).ColumnWidth = 13.73

# --- sheet selections --------------------------------------------------
$wsAge.Range("H13").Select() | Out-Null

# Activating location_id both moves tabSelected from prevalence_c429 to
# location_id and updates the workbook's bookViews/activeTab index.
$wsLoc.Activate()
$wsLoc.Range("K9").Select() | Out-Null
